$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = New-Object 'object[,]' 24,1
$col[0,0] = 0.7689166990269314
$col[1,0] = 0.6903851675846511
$col[2,0] = 0.6422450569735076
$col[3,0] = 0.6226483097021571
$col[4,0] = 0.6193955663076167
$col[5,0] = 0.6419806831131325
$col[6,0] = 0.7418232844613044
$col[7,0] = 0.9382044125550806
$col[8,0] = 1.082815134587236
$col[9,0] = 1.148668347770524
$col[10,0] = 1.173614386696499
$col[11,0] = 1.168241434091556
$col[12,0] = 1.150720502067827
$col[13,0] = 1.139989546755999
$col[14,0] = 1.078512770076145
$col[15,0] = 1.040815683081121
$col[16,0] = 1.019139881086062
$col[17,0] = 1.011801987410195
$col[18,0] = 1.044827930653469
$col[19,0] = 1.155866592722418
$col[20,0] = 1.228487900970435
$col[21,0] = 1.189724219088191
$col[22,0] = 1.04301400355132
$col[23,0] = 0.8850180759411046
$ws.Range("B2:B25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 0.1311608951239407
$col[1,0] = 0.1144111504106888
$col[2,0] = 0.1040893375892722
$col[3,0] = 0.09987393089295438
$col[4,0] = 0.09917341763946297
$col[5,0] = 0.1040325240331867
$col[6,0] = 0.125393473237807
$col[7,0] = 0.1669777327548729
$col[8,0] = 0.1973365775169782
$col[9,0] = 0.2111042668061316
$col[10,0] = 0.2163114081719471
$col[11,0] = 0.2151902448979399
$col[12,0] = 0.2115327907910967
$col[13,0] = 0.2092916559783475
$col[14,0] = 0.1964359465617065
$col[15,0] = 0.1885382749730411
$col[16,0] = 0.1839917383981344
$col[17,0] = 0.1824516788385893
$col[18,0] = 0.1893794118395249
$col[19,0] = 0.2126072489095066
$col[20,0] = 0.2277506613019682
$col[21,0] = 0.2196718322737468
$col[22,0] = 0.1889991527210952
$col[23,0] = 0.155761463523163
$ws.Range("C2:C25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 0.2617112972241102
$col[1,0] = 0.2593761002322026
$col[2,0] = 0.2580374458545549
$col[3,0] = 0.2575159405166971
$col[4,0] = 0.2574307971951271
$col[5,0] = 0.2580303153444703
$col[6,0] = 0.2608864144692973
$col[7,0] = 0.2672394484871177
$col[8,0] = 0.2723625055370604
$col[9,0] = 0.2747914141814363
$col[10,0] = 0.2757252629743192
$col[11,0] = 0.2755235172097628
$col[12,0] = 0.2748679608613003
$col[13,0] = 0.2744682442231436
$col[14,0] = 0.2722057442371408
$col[15,0] = 0.270842922910191
$col[16,0] = 0.270068330976045
$col[17,0] = 0.2698076613841351
$col[18,0] = 0.27098703903863
$col[19,0] = 0.2750601323024
$col[20,0] = 0.2778041316777404
$col[21,0] = 0.276332130061931
$col[22,0] = 0.2709218563715439
$col[23,0] = 0.2654405738057335
$ws.Range("D2:D25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 1.296422861156174
$col[1,0] = 1.293609930073984
$col[2,0] = 1.292673066549746
$col[3,0] = 1.292489912473513
$col[4,0] = 1.292471493507591
$col[5,0] = 1.292669792376032
$col[6,0] = 1.295288874020557
$col[7,0] = 1.306702159912874
$col[8,0] = 1.318926802277872
$col[9,0] = 1.325324683291456
$col[10,0] = 1.327867904013175
$col[11,0] = 1.327314815117475
$col[12,0] = 1.325531500110543
$col[13,0] = 1.324454862851439
$col[14,0] = 1.31852553742597
$col[15,0] = 1.315102525122413
$col[16,0] = 1.313212456588914
$col[17,0] = 1.312586034079487
$col[18,0] = 1.315458758441778
$col[19,0] = 1.326052031938232
$col[20,0] = 1.333677678127501
$col[21,0] = 1.329543415534616
$col[22,0] = 1.315297462816417
$col[23,0] = 1.302941245980435
$ws.Range("F2:F25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 0.002434700422176341
$col[1,0] = 0.002437619918592648
$col[2,0] = 0.00243950914249302
$col[3,0] = 0.0024403033911379
$col[4,0] = 0.002440436749862078
$col[5,0] = 0.002439519755078482
$col[6,0] = 0.002435687050749556
$col[7,0] = 0.00242893458446605
$col[8,0] = 0.002424434277044085
$col[9,0] = 0.002422486015120335
$col[10,0] = 0.00242176241268391
$col[11,0] = 0.002421917624554359
$col[12,0] = 0.002422426200602394
$col[13,0] = 0.002422739560289404
$col[14,0] = 0.002424563585817498
$col[15,0] = 0.002425707862211929
$col[16,0] = 0.002426375337498775
$col[17,0] = 0.002426602935357704
$col[18,0] = 0.002425585088148986
$col[19,0] = 0.002422276435733599
$col[20,0] = 0.002420196555477094
$col[21,0] = 0.002421299098499189
$col[22,0] = 0.002425640564310441
$col[23,0] = 0.002430680059626191
$ws.Range("G2:G25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 0.2675860458324877
$col[1,0] = 0.2675714652564523
$col[2,0] = 0.2677018638457582
$col[3,0] = 0.2677900350991607
$col[4,0] = 0.2678067915298143
$col[5,0] = 0.2677029111186897
$col[6,0] = 0.2675520904239406
$col[7,0] = 0.2683627387609349
$col[8,0] = 0.2696343163801203
$col[9,0] = 0.2703599349861889
$col[10,0] = 0.2706558913652231
$col[11,0] = 0.2705912095856604
$col[12,0] = 0.2703838589103
$col[13,0] = 0.2702596094575895
$col[14,0] = 0.2695898584392467
$col[15,0] = 0.2692166942999279
$col[16,0] = 0.2690159116021178
$col[17,0] = 0.2689503085055307
$col[18,0] = 0.2692549845625649
$col[19,0] = 0.270444187862914
$col[20,0] = 0.2713448644252878
$col[21,0] = 0.2708528537390862
$col[22,0] = 0.2692376306929276
$col[23,0] = 0.2680248340670663
$ws.Range("J2:J25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 0.3942602305829652
$col[1,0] = 0.3710925341965137
$col[2,0] = 0.3570038032193494
$col[3,0] = 0.3512970939325015
$col[4,0] = 0.3503515950524658
$col[5,0] = 0.3569267002233119
$col[6,0] = 0.3862438752183124
$col[7,0] = 0.444807254720395
$col[8,0] = 0.4884798973050977
$col[9,0] = 0.5084866790282589
$col[10,0] = 0.5160826232032036
$col[11,0] = 0.5144458244911903
$col[12,0] = 0.509111206642352
$col[13,0] = 0.5058461700250376
$col[14,0] = 0.4871751957563504
$col[15,0] = 0.4757567672697007
$col[16,0] = 0.4692023732943227
$col[17,0] = 0.4669854417673989
$col[18,0] = 0.4769709161765689
$col[19,0] = 0.5106775779969226
$col[20,0] = 0.5328221579237891
$col[21,0] = 0.5209927266811434
$col[22,0] = 0.4764219676224215
$col[23,0] = 0.4288502889005699
$ws.Range("M2:M25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 1.23369419561373
$col[1,0] = 1.246447606920917
$col[2,0] = 1.254779869254303
$col[3,0] = 1.258301497817154
$col[4,0] = 1.258893882473092
$col[5,0] = 1.254826852296862
$col[6,0] = 1.237987478314594
$col[7,0] = 1.208944245094216
$col[8,0] = 1.190029701217163
$col[9,0] = 1.181950886138424
$col[10,0] = 1.178967215261281
$col[11,0] = 1.179606439202296
$col[12,0] = 1.181703902371353
$col[13,0] = 1.182998505700553
$col[14,0] = 1.190568248101862
$col[15,0] = 1.195346657245985
$col[16,0] = 1.198144532563767
$col[17,0] = 1.199100339049515
$col[18,0] = 1.194832868197501
$col[19,0] = 1.181085774511182
$col[20,0] = 1.172541933896376
$col[21,0] = 1.177061607802926
$col[22,0] = 1.195064994284273
$col[23,0] = 1.216375485960612
$ws.Range("N2:N25").Value = $col

$col = New-Object 'object[,]' 24,1
$col[0,0] = 2.864332284145007
$col[1,0] = 2.862141862779936
$col[2,0] = 2.862699357097938
$col[3,0] = 2.863404395476664
$col[4,0] = 2.863550312527764
$col[5,0] = 2.862706931377431
$col[6,0] = 2.863181853830866
$col[7,0] = 2.879236275420141
$col[8,0] = 2.900296591568605
$col[9,0] = 2.911899577882963
$col[10,0] = 2.916584875137687
$col[11,0] = 2.915562837596156
$col[12,0] = 2.91227919427871
$col[13,0] = 2.910305849144009
$col[14,0] = 2.899579064201333
$col[15,0] = 2.893517019891533
$col[16,0] = 2.890220624267386
$col[17,0] = 2.88913719096422
$col[18,0] = 2.894142631434676
$col[19,0] = 2.913235763492622
$col[20,0] = 2.92741360875516
$col[21,0] = 2.919690914988848
$col[22,0] = 2.893859204538785
$col[23,0] = 2.873269375179973
$ws.Range("O2:O25").Value = $col

Write-Output "applied changes"